# Update the "Förändrad" (changed) date in column C for every data row
# (rows 2-244) from 2023-09-15 (serial 45184) to 2023-09-17 (serial 45186).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 244; $r++) {
    $ws.Cells.Item($r, 3).Value = 45186
}

# Add a friendly display-text second argument to the HYPERLINK() formulas
# in columns S, T, U, V, W, X, Y (rows 2-21) - the text is the record's
# "Beteckning" (column A) value, matching the target filename used in the
# link's URL.
$rowsWithLinks = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21)
$rowsWithU     = @(2,3,4,6,8,11,13,14,15,21)

foreach ($r in $rowsWithLinks) {
    $bet = $ws.Cells.Item($r, 1).Text

    $ws.Range("S$r").Formula = '=HYPERLINK("https://klasma.github.io/Logging_NORBERG/artfynd/' + $bet + '.xlsx", "' + $bet + '")'
    $ws.Range("T$r").Formula = '=HYPERLINK("https://klasma.github.io/Logging_NORBERG/kartor/' + $bet + '.png", "' + $bet + '")'

    if ($rowsWithU -contains $r) {
        $ws.Range("U$r").Formula = '=HYPERLINK("https://klasma.github.io/Logging_NORBERG/knärot/' + $bet + '.png", "' + $bet + '")'
    }

    $ws.Range("V$r").Formula = '=HYPERLINK("https://klasma.github.io/Logging_NORBERG/klagomål/' + $bet + '.docx", "' + $bet + '")'
    $ws.Range("W$r").Formula = '=HYPERLINK("https://klasma.github.io/Logging_NORBERG/klagomålsmail/' + $bet + '.docx", "' + $bet + '")'
    $ws.Range("X$r").Formula = '=HYPERLINK("https://klasma.github.io/Logging_NORBERG/tillsyn/' + $bet + '.docx", "' + $bet + '")'
    $ws.Range("Y$r").Formula = '=HYPERLINK("https://klasma.github.io/Logging_NORBERG/tillsynsmail/' + $bet + '.docx", "' + $bet + '")'
}
